# Auto-update of Leve price/profit tracking cells (scheduled market-price refresh).
# Values come from an external price feed; this mirrors that refreshed snapshot
# onto the per-job "Profits" sheets (one sheet per crafting discipline).

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value2 = 234.54546
$ws.Range("I6").Value2 = 234.54546
$ws.Range("K6").Value2 = 703.6363799999999
$ws.Range("M6").Value2 = -591.6363799999999
$ws.Range("H31").Value2 = 170
$ws.Range("I31").Value2 = 170
$ws.Range("K31").Value2 = 510
$ws.Range("M31").Value2 = -280
$ws.Range("H40").Value2 = 2223.1538
$ws.Range("I40").Value2 = 2250.1
$ws.Range("J40").Value2 = 2133.3333
$ws.Range("K40").Value2 = 2250.1
$ws.Range("L40").Value2 = 2133.3333
$ws.Range("M40").Value2 = -2075.1
$ws.Range("N40").Value2 = -2483.3333
$ws.Range("H46").Value2 = 1416.6666
$ws.Range("I46").Value2 = 0
$ws.Range("J46").Value2 = 1416.6666
$ws.Range("K46").Value2 = 0
$ws.Range("L46").Value2 = 4249.9998
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value2 = -4487.9998
$ws.Range("H60").Value2 = 1416.6666
$ws.Range("I60").Value2 = 0
$ws.Range("J60").Value2 = 1416.6666
$ws.Range("K60").Value2 = 0
$ws.Range("L60").Value2 = 4249.9998
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value2 = -5217.9998
$ws.Range("H64").Value2 = 4750
$ws.Range("I64").Value2 = 0
$ws.Range("J64").Value2 = 4750
$ws.Range("K64").Value2 = 0
$ws.Range("L64").Value2 = 4750
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value2 = -5246
$ws.Range("H67").Value2 = 4750
$ws.Range("I67").Value2 = 0
$ws.Range("J67").Value2 = 4750
$ws.Range("K67").Value2 = 0
$ws.Range("L67").Value2 = 4750
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value2 = -6466
$ws.Range("H68").Value2 = 28647.5
$ws.Range("J68").Value2 = 28647.5
$ws.Range("L68").Value2 = 28647.5
$ws.Range("N68").Value2 = -30145.5
$ws.Range("H71").Value2 = 28647.5
$ws.Range("J71").Value2 = 28647.5
$ws.Range("L71").Value2 = 85942.5
$ws.Range("N71").Value2 = -93430.5
$ws.Range("H76").Value2 = 2959.875
$ws.Range("I76").Value2 = 2720.6155
$ws.Range("J76").Value2 = 3996.6667
$ws.Range("K76").Value2 = 2720.6155
$ws.Range("L76").Value2 = 3996.6667
$ws.Range("M76").Value2 = -2405.6155
$ws.Range("N76").Value2 = -4626.6667
$ws.Range("H79").Value2 = 2959.875
$ws.Range("I79").Value2 = 2720.6155
$ws.Range("J79").Value2 = 3996.6667
$ws.Range("K79").Value2 = 2720.6155
$ws.Range("L79").Value2 = 3996.6667
$ws.Range("M79").Value2 = -1628.6155
$ws.Range("N79").Value2 = -6180.6667
$ws.Range("H86").Value2 = 2522.889
$ws.Range("I86").Value2 = 1940.4
$ws.Range("J86").Value2 = 3251
$ws.Range("K86").Value2 = 1940.4
$ws.Range("L86").Value2 = 3251
$ws.Range("M86").Value2 = -817.4000000000001
$ws.Range("N86").Value2 = -5497
$ws.Range("H89").Value2 = 2522.889
$ws.Range("I89").Value2 = 1940.4
$ws.Range("J89").Value2 = 3251
$ws.Range("K89").Value2 = 9702
$ws.Range("L89").Value2 = 16255
$ws.Range("M89").Value2 = -4086
$ws.Range("N89").Value2 = -27487
$ws.Range("H98").Value2 = 44774.39
$ws.Range("I98").Value2 = 65459.582
$ws.Range("J98").Value2 = 2025
$ws.Range("K98").Value2 = 65459.582
$ws.Range("L98").Value2 = 2025
$ws.Range("M98").Value2 = -63961.582
$ws.Range("N98").Value2 = -5021
$ws.Range("H122").Value2 = 44774.39
$ws.Range("I122").Value2 = 65459.582
$ws.Range("J122").Value2 = 2025
$ws.Range("K122").Value2 = 196378.746
$ws.Range("L122").Value2 = 6075
$ws.Range("M122").Value2 = -193928.746
$ws.Range("N122").Value2 = -10975
$ws.Range("H131").Value2 = 863.3333
$ws.Range("I131").Value2 = 496.15384
$ws.Range("K131").Value2 = 1488.46152
$ws.Range("M131").Value2 = 3551.53848
$ws.Range("H132").Value2 = 1521102.2
$ws.Range("I132").Value2 = 2041406
$ws.Range("K132").Value2 = 6124218
$ws.Range("M132").Value2 = -6121688

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value2 = 1117102.9
$ws.Range("I33").Value2 = 2002485.2
$ws.Range("J33").Value2 = 10375
$ws.Range("K33").Value2 = 2002485.2
$ws.Range("L33").Value2 = 10375
$ws.Range("M33").Value2 = -2002156.2
$ws.Range("N33").Value2 = -11033
$ws.Range("H63").Value2 = 2046.6666
$ws.Range("I63").Value2 = 1912.5
$ws.Range("J63").Value2 = 2200
$ws.Range("K63").Value2 = 1912.5
$ws.Range("L63").Value2 = 2200
$ws.Range("M63").Value2 = -1226.5
$ws.Range("N63").Value2 = -3572
$ws.Range("H66").Value2 = 2046.6666
$ws.Range("I66").Value2 = 1912.5
$ws.Range("J66").Value2 = 2200
$ws.Range("K66").Value2 = 9562.5
$ws.Range("L66").Value2 = 11000
$ws.Range("M66").Value2 = -6130.5
$ws.Range("N66").Value2 = -17864
$ws.Range("H74").Value2 = 13159077
$ws.Range("I74").Value2 = 15152528
$ws.Range("J74").Value2 = 2300
$ws.Range("K74").Value2 = 15152528
$ws.Range("L74").Value2 = 2300
$ws.Range("M74").Value2 = -15151654
$ws.Range("N74").Value2 = -4048
$ws.Range("H77").Value2 = 13159077
$ws.Range("I77").Value2 = 15152528
$ws.Range("J77").Value2 = 2300
$ws.Range("K77").Value2 = 75762640
$ws.Range("L77").Value2 = 11500
$ws.Range("M77").Value2 = -75758272
$ws.Range("N77").Value2 = -20236

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value2 = 0
$ws.Range("I26").Value2 = 0
$ws.Range("J26").Value2 = 0
$ws.Range("K26").Value2 = 0
$ws.Range("L26").Value2 = 0
$ws.Range("M26").ClearContents()
$ws.Range("N26").ClearContents()
$ws.Range("H76").Value2 = 47500
$ws.Range("J76").Value2 = 47500
$ws.Range("L76").Value2 = 47500
$ws.Range("N76").Value2 = -48130
$ws.Range("H79").Value2 = 47500
$ws.Range("J79").Value2 = 47500
$ws.Range("L79").Value2 = 47500
$ws.Range("N79").Value2 = -49684
$ws.Range("H86").Value2 = 1697.381
$ws.Range("I86").Value2 = 1554.6875
$ws.Range("K86").Value2 = 1554.6875
$ws.Range("M86").Value2 = -431.6875
$ws.Range("H89").Value2 = 1697.381
$ws.Range("I89").Value2 = 1554.6875
$ws.Range("K89").Value2 = 7773.4375
$ws.Range("M89").Value2 = -2157.4375
$ws.Range("H105").Value2 = 2727.7693
$ws.Range("I105").Value2 = 2308.3333
$ws.Range("J105").Value2 = 3087.2856
$ws.Range("K105").Value2 = 2308.3333
$ws.Range("L105").Value2 = 3087.2856
$ws.Range("M105").Value2 = -561.3332999999998
$ws.Range("N105").Value2 = -6581.2856

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value2 = 2048.3635
$ws.Range("I13").Value2 = 1670.2222
$ws.Range("J13").Value2 = 3750
$ws.Range("K13").Value2 = 1670.2222
$ws.Range("L13").Value2 = 3750
$ws.Range("M13").Value2 = -1531.2222
$ws.Range("N13").Value2 = -4028
$ws.Range("H62").Value2 = 202400
$ws.Range("J62").Value2 = 4000
$ws.Range("L62").Value2 = 4000
$ws.Range("N62").Value2 = -5248
$ws.Range("H65").Value2 = 202400
$ws.Range("J65").Value2 = 4000
$ws.Range("L65").Value2 = 20000
$ws.Range("N65").Value2 = -26240
$ws.Range("H88").Value2 = 0
$ws.Range("J88").Value2 = 0
$ws.Range("L88").Value2 = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value2 = 0
$ws.Range("J91").Value2 = 0
$ws.Range("L91").Value2 = 0
$ws.Range("N91").ClearContents()
$ws.Range("H122").Value2 = 4567.5312
$ws.Range("I122").Value2 = 5182.88
$ws.Range("K122").Value2 = 15548.64
$ws.Range("M122").Value2 = -13098.64
$ws.Range("H134").Value2 = 20834676
$ws.Range("I134").Value2 = 33334378
$ws.Range("J134").Value2 = 1839.5555
$ws.Range("K134").Value2 = 100003134
$ws.Range("L134").Value2 = 5518.666499999999
$ws.Range("M134").Value2 = -100000599
$ws.Range("N134").Value2 = -10588.6665

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value2 = 4162.6
$ws.Range("I82").Value2 = 1013
$ws.Range("J82").Value2 = 4950
$ws.Range("K82").Value2 = 3039
$ws.Range("L82").Value2 = 14850
$ws.Range("M82").Value2 = -2633
$ws.Range("N82").Value2 = -15662
$ws.Range("H85").Value2 = 4162.6
$ws.Range("I85").Value2 = 1013
$ws.Range("J85").Value2 = 4950
$ws.Range("K85").Value2 = 3039
$ws.Range("L85").Value2 = 14850
$ws.Range("M85").Value2 = -1635
$ws.Range("N85").Value2 = -17658
$ws.Range("H88").Value2 = 4128
$ws.Range("J88").Value2 = 4128
$ws.Range("L88").Value2 = 12384
$ws.Range("N88").Value2 = -13240
$ws.Range("H91").Value2 = 4128
$ws.Range("J91").Value2 = 4128
$ws.Range("L91").Value2 = 12384
$ws.Range("N91").Value2 = -15348
$ws.Range("H94").Value2 = 5000
$ws.Range("I94").Value2 = 5000
$ws.Range("K94").Value2 = 15000
$ws.Range("M94").Value2 = -14324

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value2 = 892
$ws.Range("I31").Value2 = 892
$ws.Range("K31").Value2 = 892
$ws.Range("M31").Value2 = -600
$ws.Range("H37").Value2 = 892
$ws.Range("I37").Value2 = 892
$ws.Range("K37").Value2 = 892
$ws.Range("M37").Value2 = -615
$ws.Range("H70").Value2 = 4137.6875
$ws.Range("I70").Value2 = 4037.4
$ws.Range("K70").Value2 = 4037.4
$ws.Range("M70").Value2 = -3767.4
$ws.Range("H73").Value2 = 4137.6875
$ws.Range("I73").Value2 = 4037.4
$ws.Range("K73").Value2 = 4037.4
$ws.Range("M73").Value2 = -3101.4
$ws.Range("H80").Value2 = 145115.14
$ws.Range("I80").Value2 = 2250
$ws.Range("J80").Value2 = 202261.2
$ws.Range("K80").Value2 = 2250
$ws.Range("L80").Value2 = 202261.2
$ws.Range("M80").Value2 = -1252
$ws.Range("N80").Value2 = -204257.2
$ws.Range("H82").Value2 = 33000
$ws.Range("H83").Value2 = 145115.14
$ws.Range("I83").Value2 = 2250
$ws.Range("J83").Value2 = 202261.2
$ws.Range("K83").Value2 = 11250
$ws.Range("L83").Value2 = 1011306
$ws.Range("M83").Value2 = -6258
$ws.Range("N83").Value2 = -1021290
$ws.Range("H85").Value2 = 33000
$ws.Range("H88").Value2 = 0
$ws.Range("J88").Value2 = 0
$ws.Range("L88").Value2 = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value2 = 0
$ws.Range("J91").Value2 = 0
$ws.Range("L91").Value2 = 0
$ws.Range("N91").ClearContents()

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value2 = 1256.5625
$ws.Range("I122").Value2 = 1061.1
$ws.Range("J122").Value2 = 1582.3334
$ws.Range("K122").Value2 = 3183.3
$ws.Range("L122").Value2 = 4747.0002
$ws.Range("M122").Value2 = -733.2999999999997
$ws.Range("N122").Value2 = -9647.0002
